# "contingencies with rene fine"
#
# Populate Sheet1 with a tiny 2x2 "disconnected elements" block:
#   B1 = 0                         (bold, thin box border, centered/top aligned)
#   A2 = 0                         (bold, thin box border, centered/top aligned)
#   B2 = "disconnected_elements"   (plain text, default formatting)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlExcel constants used below (no Excel type library available, so use literals):
#   xlContinuous      =  1   (Borders.LineStyle)
#   xlCenter          = -4108 (HorizontalAlignment)
#   xlTop             = -4160 (VerticalAlignment)
#   xlPasteFormats    = -4122 (PasteSpecial Paste:=xlPasteFormats)

$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the bold / boxed / centered-top style on B1 first ...
$rng = $ws.Range("B1")
$rng.Borders.LineStyle = 1
$rng.Font.Bold = $true
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4160

# ... then clone just the formatting onto A2 so both cells share one cell style
# (setting the properties on A2 directly would otherwise mint a second, slightly
# different style record).
$rng.Copy()
$ws.Range("A2").PasteSpecial(-4122)
